$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/string formatting so numeric-looking values (with significant
# trailing zeros) and percentages are preserved as text rather than being
# auto-converted to numbers by Excel.
$cells = @{
    "D2" = "301.27"
    "E2" = "0.88%"
    "D3" = "31.80"
    "E3" = "1.71%"
    "D4" = "5.095"
    "E4" = "-0.59%"
    "D5" = "0.07808"
    "E5" = "-2.85%"
    "D6" = "2.248"
    "E6" = "-13.77%"
    "D7" = "7.790"
    "E7" = "-0.60%"
    "D8" = "3.815"
    "E8" = "-0.24%"
    "D9" = "0.9172"
    "E9" = "-0.12%"
    "D10" = "0.1761"
    "E10" = "1.40%"
    "D11" = "0.07538"
    "E11" = "3.14%"
    "E12" = "7.78%"
    "D13" = "0.03032"
    "E13" = "0.04%"
    "D14" = "0.1003"
    "E14" = "0.60%"
    "D15" = "0.001504"
    "E15" = "0.04%"
    "D16" = "0.005891"
    "E16" = "-1.35%"
    "D17" = "3.469"
    "E17" = "-0.89%"
    "E18" = "0.00%"
    "E19" = "0.24%"
    "E20" = "-0.10%"
    "D21" = "4.224"
    "E21" = "-8.92%"
    "D22" = "0.1816"
    "E22" = "13.62%"
    "D23" = "0.04595"
    "E23" = "0.82%"
    "E24" = "-0.90%"
    "D25" = "0.004470"
    "E25" = "0.48%"
    "D26" = "0.0001248"
    "E26" = "5.80%"
    "E27" = "-1.44%"
    "D39" = "0.01776"
    "E39" = "-3.32%"
    "E40" = "5.99%"
    "D41" = "0.007378"
    "E41" = "5.33%"
    "E42" = "1.30%"
    "D43" = "0.002187"
    "E43" = "-2.35%"
    "D44" = "0.01027"
    "E44" = "4.53%"
    "D45" = "0.00006279"
    "E45" = "-3.19%"
    "E46" = "-0.17%"
    "D48" = "0.7369"
    "E48" = "-10.19%"
    "D49" = "0.00002097"
    "E49" = "-0.17%"
    "D50" = "0.0001998"
    "E50" = "-0.17%"
}

foreach ($ref in $cells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cells[$ref]
    $cell.Style = "Normal"
}
